$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings are not
# coerced into numbers (matches original inline-string cell content).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.569.11'
$ws.Range('E2').Value = '  +4.84%  '
$ws.Range('D3').Value = '2.751.21'
$ws.Range('E3').Value = '  +4.78%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '587.73'
$ws.Range('E5').Value = '  -2.86%  '
$ws.Range('D6').Value = '152.55'
$ws.Range('E6').Value = '  +5.58%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '0.607'
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('D9').Value = '2.751.71'
$ws.Range('E9').Value = '  +4.22%  '
$ws.Range('D10').Value = '6.71'
$ws.Range('E10').Value = '  +2.54%  '
$ws.Range('D11').Value = '0.112'
$ws.Range('E11').Value = '  +5.81%  '
$ws.Range('E12').Value = '  +4.00%  '
$ws.Range('D13').Value = '0.160'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').Value = '3.213.30'
$ws.Range('E14').Value = '  +3.93%  '
$ws.Range('D15').Value = '26.34'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '63.408.38'
$ws.Range('E16').Value = '  +4.57%  '
$ws.Range('E17').Value = '  +7.92%  '
$ws.Range('D18').Value = '2.743.14'
$ws.Range('E18').Value = '  +4.00%  '
$ws.Range('D19').Value = '12.05'
$ws.Range('E19').Value = '  +4.51%  '
$ws.Range('E20').Value = '  +3.19%  '
$ws.Range('D21').Value = '362.09'
$ws.Range('E21').Value = '  +3.73%  '
$ws.Range('D22').Value = '7.02'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').Value = '0.538'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').Value = '65.60'
$ws.Range('E25').Value = '  +2.98%  '
$ws.Range('E26').Value = '  +4.23%  '
$ws.Range('E27').Value = '  +6.84%  '
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('D29').Value = '0.0₃0895'
$ws.Range('E29').Value = '  +11.36%  '
$ws.Range('D30').Value = '2.04'
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('D31').Value = '7.10'
$ws.Range('E31').Value = '  +8.70%  '
$ws.Range('D32').Value = '173.55'
$ws.Range('E32').Value = '  +2.61%  '
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').Value = '  +18.23%  '
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').Value = '20.54'
$ws.Range('E35').Value = '  +4.86%  '
$ws.Range('D36').Value = '4.77'
$ws.Range('E36').Value = '  +8.05%  '
$ws.Range('D37').Value = '1.44'
$ws.Range('E37').Value = '  +9.92%  '
$ws.Range('D38').Value = '1.80'
$ws.Range('E38').Value = '  +10.11%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +16.99%  '
$ws.Range('D40').Value = '344.45'
$ws.Range('E40').Value = '  +4.34%  '
$ws.Range('D41').Value = '4.24'
$ws.Range('E41').Value = '  +5.49%  '
$ws.Range('D42').Value = '38.93'
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('D43').Value = '5.55'
$ws.Range('E43').Value = '  +6.18%  '
$ws.Range('D44').Value = '21.85'
$ws.Range('E44').Value = '  +8.38%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '22.09'
$ws.Range('E45').Value = '  +9.50%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '141.52'
$ws.Range('E46').Value = '  +5.24%  '
$ws.Range('D47').Value = '0.0589'
$ws.Range('E47').Value = '  +5.78%  '
$ws.Range('D48').Value = '0.643'
$ws.Range('E48').Value = '  +5.41%  '
$ws.Range('D49').Value = '0.0256'
$ws.Range('E49').Value = '  +5.46%  '
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  -0.21%  '
